$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 32 (the most recent "Primera" entry for 2022-02-15) twice,
# inserting the copies right below it. This pushes the old "Segunda" row
# (previously row 33) down to row 35, and leaves two fresh copies of the
# old row 32 data at rows 33 and 34.
$ws.Rows("32:32").Copy()
$ws.Rows("33:33").Insert()
$ws.Rows("32:32").Copy()
$ws.Rows("33:33").Insert()

# Row 32 becomes the new weekly "Especial" quality entry for 2023-02-27.
$ws.Cells.Item(32, 4).Value2 = 44984
$ws.Cells.Item(32, 12).Value2 = "Especial"
$ws.Cells.Item(32, 13).Value2 = 160
$ws.Cells.Item(32, 14).Value2 = 13000
$ws.Cells.Item(32, 15).Value2 = 14000
$ws.Cells.Item(32, 16).Value2 = 13500
$ws.Cells.Item(32, 19).Value2 = 750

# Row 33 keeps the old row-32 "Primera" values but is dated with the new
# week (2023-02-27).
$ws.Cells.Item(33, 4).Value2 = 44984

# Rows 34 and 35 retain the original data (old rows 32 and 33), unchanged.
